{"js": "// Replace each two-digit multiplication expression in the document with\n// its new value. Every occurrence is unique in the document (no string\n// appears twice among the \"before\" set, and none of the \"after\" values\n// collide with any \"before\" value), so a direct search-and-replace per\n// pair is safe regardless of processing order.\nconst replacements = [\n  [\"43\u00d770=\", \"21\u00d726=\"],\n  [\"37\u00d734=\", \"65\u00d799=\"],\n  [\"53\u00d719=\", \"30\u00d726=\"],\n  [\"81\u00d779=\", \"21\u00d793=\"],\n  [\"41\u00d781=\", \"73\u00d780=\"],\n  [\"33\u00d720=\", \"34\u00d722=\"],\n  [\"70\u00d741=\", \"95\u00d734=\"],\n  [\"56\u00d783=\", \"48\u00d759=\"],\n  [\"97\u00d723=\", \"66\u00d742=\"],\n  [\"20\u00d722=\", \"22\u00d743=\"],\n  [\"23\u00d793=\", \"39\u00d799=\"],\n  [\"90\u00d774=\", \"92\u00d733=\"],\n  [\"16\u00d795=\", \"95\u00d747=\"],\n  [\"91\u00d787=\", \"28\u00d763=\"],\n  [\"77\u00d732=\", \"64\u00d745=\"],\n  [\"18\u00d766=\", \"65\u00d739=\"],\n  [\"27\u00d721=\", \"78\u00d739=\"],\n  [\"92\u00d768=\", \"90\u00d781=\"],\n  [\"24\u00d740=\", \"87\u00d740=\"],\n  [\"17\u00d785=\", \"84\u00d713=\"],\n  [\"16\u00d756=\", \"13\u00d742=\"],\n  [\"99\u00d782=\", \"93\u00d752=\"],\n  [\"85\u00d734=\", \"20\u00d769=\"],\n  [\"24\u00d737=\", \"90\u00d755=\"],\n  [\"44\u00d716=\", \"54\u00d787=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the document with\n# its new value. Every occurrence is unique in the document (no string\n# appears twice among the \"before\" set, and none of the \"after\" values\n# collide with any \"before\" value), so a direct Find/Replace per pair is\n# safe regardless of processing order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"43\u00d770=\", \"21\u00d726=\"),\n    @(\"37\u00d734=\", \"65\u00d799=\"),\n    @(\"53\u00d719=\", \"30\u00d726=\"),\n    @(\"81\u00d779=\", \"21\u00d793=\"),\n    @(\"41\u00d781=\", \"73\u00d780=\"),\n    @(\"33\u00d720=\", \"34\u00d722=\"),\n    @(\"70\u00d741=\", \"95\u00d734=\"),\n    @(\"56\u00d783=\", \"48\u00d759=\"),\n    @(\"97\u00d723=\", \"66\u00d742=\"),\n    @(\"20\u00d722=\", \"22\u00d743=\"),\n    @(\"23\u00d793=\", \"39\u00d799=\"),\n    @(\"90\u00d774=\", \"92\u00d733=\"),\n    @(\"16\u00d795=\", \"95\u00d747=\"),\n    @(\"91\u00d787=\", \"28\u00d763=\"),\n    @(\"77\u00d732=\", \"64\u00d745=\"),\n    @(\"18\u00d766=\", \"65\u00d739=\"),\n    @(\"27\u00d721=\", \"78\u00d739=\"),\n    @(\"92\u00d768=\", \"90\u00d781=\"),\n    @(\"24\u00d740=\", \"87\u00d740=\"),\n    @(\"17\u00d785=\", \"84\u00d713=\"),\n    @(\"16\u00d756=\", \"13\u00d742=\"),\n    @(\"99\u00d782=\", \"93\u00d752=\"),\n    @(\"85\u00d734=\", \"20\u00d769=\"),\n    @(\"24\u00d737=\", \"90\u00d755=\"),\n    @(\"44\u00d716=\", \"54\u00d787=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute(\n        $oldText,       # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap: wdFindContinue\n        $false,         # Format\n        $newText,       # ReplaceWith\n        2               # Replace: wdReplaceAll\n    ) | Out-Null\n}\n"}
